# Financial Statement Output - re-label/re-point rows 2-37 to the
# "cash flows" section of the statement (the author is still untangling
# which line items belong where). Row 38 onward is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row([int]$r, [string]$label, $b = $null, $c = $null, $d = $null) {
    $ws.Cells.Item($r, 1).Value = $label
    if ($null -ne $b) {
        $ws.Cells.Item($r, 2).Value = $b
        $ws.Cells.Item($r, 3).Value = $c
        $ws.Cells.Item($r, 4).Value = $d
    } else {
        $ws.Range($ws.Cells.Item($r, 2), $ws.Cells.Item($r, 4)).ClearContents()
    }
}

Set-Row 2  "cash flows from operating activities"
Set-Row 3  "adjustments to reconcile net income (loss) to net cash provided by operating activities:"
Set-Row 4  "depreciation amortization and impairment" 2322 2154 2154
Set-Row 5  "stock-based compensation" 1734 898 898
Set-Row 6  "amortization of debt discounts and issuance costs" 180 188 188
Set-Row 7  "inventory and purchase commitments write-downs" 202 193 193
Set-Row 8  "loss on disposals of fixed assets" 117 146 146
Set-Row 9  "foreign currency transaction net loss (gain)" 114 0 0
Set-Row 10 "non-cash interest and other operating activities" 228 186 186
Set-Row 11 "changes in operating assets and liabilities net of effect of business combinations:"
Set-Row 12 "prepaid expenses and other current assets" 0 0 0
Set-Row 13 "accounts payable and accrued liabilities" 2102 646 646
Set-Row 14 "deferred revenue" 321 801 801
Set-Row 15 "customer deposits" 7 0 0
Set-Row 16 "other long-term liabilities" 495 0 0
Set-Row 17 "net cash provided by operating activities" 5.943 2.405 2.405
Set-Row 18 "cash flows from investing activities"
Set-Row 19 "purchases of property and equipment excluding finance leases net of sales 3157) (1327) (2101)"
Set-Row 20 "receipt of government grants" 123 46 46
Set-Row 21 "purchase of intangible assets" 0 0 0
Set-Row 22 "cash flows from financing activities"
Set-Row 23 "proceeds from issuances of convertible and other debt" 9713 10669 10669
Set-Row 24 "repayments of convertible and other debt (11623) (9161) (5247)"
Set-Row 25 "collateralized lease repayments" 0 0 0
Set-Row 26 "proceeds from exercises of stock options and other stock issuances" 417 263 263
Set-Row 27 "purchase of convertible note hedges" 0 0 0
Set-Row 28 "proceeds from investments by noncontrolling interests in subsidiaries" 24 279 279
Set-Row 29 "net cash provided by financing activities" 9.973 1529 1529
Set-Row 30 "effect of exchange rate changes on cash and cash equivalents and restricted cash" 334 8 8
Set-Row 31 "net increase in cash and cash equivalents and restricted cash" 13118 2506 2506
Set-Row 32 "cash and cash equivalents and restricted cash beginning of period" 6.783 4277 4277
Set-Row 33 "cash and cash equivalents and restricted cash end of period" 19.901 6.783 6.783
Set-Row 34 "supplemental non-cash investing and financing activities"
Set-Row 35 "acquisitions property and equipment included in liabilities" 1088 562 562
Set-Row 36 "supplemental disclosures"
Set-Row 37 "cash paid during the period for taxes net of the refunds accompanying notes are an integral part of these consolidated financial us statements. os" 58 54 54
